# ----------------------------------------------------------------------------
# Applies the resume-content update described by the commit:
#  - "21 years" -> "15+ years" in the professional summary
#  - Expanded FLEEM bullet (Progressive Change Campaign Committee)
#  - Rewrote + expanded Salsa Labs bullets (4 -> 6)
#  - Rewrote + expanded Praxis Project bullets (4 -> 7)
#  - Added a new training bullet under Lake Research Partners
#  - Added a new training bullet under The Feldman Group
# ----------------------------------------------------------------------------

$d = $word.ActiveDocument

function Replace-ExactText($doc, [string]$oldText, [string]$newText) {
    $r = $doc.Content
    $ok = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Replace-ExactText: could not find text: $oldText"
    }
}

function Get-ParaByText($doc, [string]$text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        $t = $t -replace "[\r\a\x07]+$", ""
        if ($t -eq $text) {
            return $doc.Paragraphs($i)
        }
    }
    return $null
}

function Insert-BulletAfter($doc, [string]$anchorText, [string]$newText) {
    $p = Get-ParaByText $doc $anchorText
    if ($p -eq $null) {
        throw "Insert-BulletAfter: anchor paragraph not found: $anchorText"
    }
    $p.Range.InsertParagraphAfter()
    $newP = $p.Next()
    $newP.Range.Text = $newText
}

# 1) Professional summary: "21 years" -> "15+ years"
Replace-ExactText $d `
    "Distinguished Polling, Research & Redistricting Professional with 21 years of expertise in survey methodology, consumer insights, voting behavior, and advanced data analysis. Proven track record in designing and implementing comprehensive research studies, managing cross-functional teams, and translating complex data into actionable intelligence. Expert in geospatial analysis, demographic segmentation, and consumer behavior modeling with experience serving major brands, organizations, and political candidates. Regular expert testimony and source on public opinion for journalists, with redistricting analysis used in court cases." `
    "Distinguished Polling, Research & Redistricting Professional with 15+ years of expertise in survey methodology, consumer insights, voting behavior, and advanced data analysis. Proven track record in designing and implementing comprehensive research studies, managing cross-functional teams, and translating complex data into actionable intelligence. Expert in geospatial analysis, demographic segmentation, and consumer behavior modeling with experience serving major brands, organizations, and political candidates. Regular expert testimony and source on public opinion for journalists, with redistricting analysis used in court cases."

# 2) FLEEM bullet (Progressive Change Campaign Committee)
Replace-ExactText $d `
    "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls" `
    "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys"

# 3) Salsa Labs bullets: rewrite the first four, then add two new ones
Replace-ExactText $d `
    "• Developed software solutions for political campaigns and advocacy groups" `
    "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously"

Replace-ExactText $d `
    "• Built web applications for voter engagement and campaign management" `
    "• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers"

Replace-ExactText $d `
    "• Integrated third-party APIs and data sources for campaign tools" `
    "• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill"

Replace-ExactText $d `
    "• Collaborated with political strategists to translate requirements into technical solutions" `
    "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs"

Insert-BulletAfter $d `
    "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs" `
    "• Collaborated with political strategists to translate geospatial requirements into technical solutions"

Insert-BulletAfter $d `
    "• Collaborated with political strategists to translate geospatial requirements into technical solutions" `
    "• Handled billions of records with millions of columns in high-performance CRM system"

# 4) Praxis Project bullets: rewrite the first four, then add three new ones
Replace-ExactText $d `
    "• Integrated technology solutions within organizational frameworks for social justice organizations" `
    "• Led technology operations for multi-million dollar organization while assisting in search for full-time CTO"

Replace-ExactText $d `
    "• Developed data management systems for community organizing efforts" `
    "• Directed all technology decisions and practices for massive multinational non-governmental organization"

Replace-ExactText $d `
    "• Provided technical training and support to nonprofit staff" `
    "• Developed comprehensive frameworks for internal and external technology audits"

Replace-ExactText $d `
    "• Built custom applications for community engagement and advocacy" `
    "• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research"

Insert-BulletAfter $d `
    "• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research" `
    "• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL"

Insert-BulletAfter $d `
    "• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL" `
    "• Managed technology infrastructure supporting community health initiatives across multiple countries"

Insert-BulletAfter $d `
    "• Managed technology infrastructure supporting community health initiatives across multiple countries" `
    "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# 5) New training bullet under Lake Research Partners
Insert-BulletAfter $d `
    "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding" `
    "• Trained staff on building Python tooling for report generation and analysis"

# 6) New training bullet under The Feldman Group
Insert-BulletAfter $d `
    "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL" `
    "• Trained staff on PHP/MySQL for data analysis and reporting systems"

Write-Output "All edits applied successfully."
